$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

$data = @(
    @(0, 12338.6, 13197, 11448, 0.2636765877405802),
    @(1, 12390.96666666667, 13509, 11390, 0.2740828275680542),
    @(2, 12886.8, 13738, 12124, 0.2677592992782593),
    @(3, 12463.9, 13179, 11175, 0.2601450045903524),
    @(4, 12143.5, 12875, 11001, 0.2592458883921305),
    @(5, 13056.23333333333, 13955, 11959, 0.2691202243169149),
    @(6, 12852.26666666667, 13486, 11385, 0.2706339359283447),
    @(7, 12250.03333333333, 13324, 9627, 0.265869148572286),
    @(8, 11975.03333333333, 13058, 9878, 0.2623656511306763),
    @(9, 12213.16666666667, 13125, 11159, 0.2620581229527791)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
